$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("December 19, 2021", $true, $false, $false, $false, $false,
               $true, 1, $false, "December 20, 2021", 2)
